$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns F,G,H,I
$ws.Range("F1").Value = "09-04-2025 Status"
$ws.Range("G1").Value = "09-04-2025 Time"
$ws.Range("H1").Value = "10-04-2025 Status"
$ws.Range("I1").Value = "10-04-2025 Time"

# Copy the style of D1:E1 (header style) onto F1:I1
$ws.Range("D1:E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("D1:E1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill the data rows 2-29 for F,G,H,I with "A" / "00:00:00" pattern
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 6).Value = "A"        # F
    $ws.Cells.Item($r, 7).Value = "00:00:00" # G
    $ws.Cells.Item($r, 8).Value = "A"        # H
    $ws.Cells.Item($r, 9).Value = "00:00:00" # I
}
